# Slide 6, "Text Placeholder 3" shape: the "Is missile wrap able: False" line
# is currently split across five runs:
#   "Is missile " | "wrap able" | ": " | "F" | "alse"
# The target state collapses it down to two runs:
#   "Is missile wrap able: " | "False"
# (the visible text is unchanged, only the run layout/formatting runs differ)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange

$fullText = $tr.Text
$lineStart = $fullText.IndexOf("Is missile ") + 1

# First run originally reads "Is missile " (11 characters) - rewrite it so it
# absorbs the text of the old 2nd/3rd runs ("wrap able" + ": ").
$firstRun = $tr.Characters($lineStart, 11)
$firstRun.Text = "Is missile wrap able: "

# Everything remaining on the line ("wrap able: False") now immediately
# follows the freshly written first run - rewrite it down to a single
# "False" run (replacing the old "wrap able", ": ", "F", "alse" runs).
$restStart = $lineStart + "Is missile wrap able: ".Length
$restRun = $tr.Characters($restStart, "wrap able: False".Length)
$restRun.Text = "False"
